$d = $word.ActiveDocument

# Split the "Dheeraj Chand" name paragraph into two paragraphs: the
# existing name line, followed by a new centered contact-info line.
$find = $d.Content.Find
$find.Execute(
    "Dheeraj Chand",            # FindText
    $false,                     # MatchCase
    $false,                     # MatchWholeWord
    $false,                     # MatchWildcards
    $false,                     # MatchSoundsLike
    $false,                     # MatchAllWordForms
    $true,                      # Forward
    1,                          # Wrap (wdFindContinue)
    $false,                     # Format
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2                           # Replace (wdReplaceAll)
)
